# issue #5: stock data output to json file
#
# Adds a new "property_category" column (value "stock") to the 股票
# (stock) sheet, right after the "total" column and before "date" /
# "legislator_name" / "legislator_id" — shifting those three columns one
# to the right. Also fixes a stray space in the THSR company name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new blank column at H (shifting old H:J -> I:K), carrying the
# existing column formatting along with it.
$ws.Range("H1:H2").Insert(-4161)

# New header + value for the inserted "property_category" column.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"

# Fix the stray mid-word space in the company name.
$ws.Range("B2").Value = "台灣高速鐵路股份有限公司"
